# PROS-13075 - CCRU - POS KPI 2020 change
#
# The underlying change is a re-sort (alphabetical) of the KPI "Benchmark 2020"
# rows' shared-string table. Because every row (2-14) in columns B/D/E/F held a
# reference into that table, re-sorting it effectively re-shuffles which KPI
# name text shows up on each row. We reproduce that end result directly by
# writing the final (alphabetically-ordered) KPI names into the right rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> final KPI text (columns B, D, E and F all mirror the same text).
$kpiByRow = @{
  2  = "CCH coolers quality"
  3  = "CCH coolers quality (Prime Pos/Max15/Merch STD/Occupancy/Lights&chilled)"
  4  = "CCH products present in Customers menu"
  5  = "CCH shelf share in Energy"
  6  = "CCH shelf share in Juice"
  7  = "CCH shelf share in SSD"
  8  = "CCH shelf share in Tea"
  9  = "CCH shelf share in Water"
  10 = "Number of CCH activation points in NARTD"
  11 = "Number of CCH cooler doors and/or equivalent in Customer coolers"
  12 = "Number of CCH displays points of interaction"
  13 = "Number of NCB core assortment available in-store"
  14 = "Number of SSD core assortment available in-store"
}

foreach ($r in 2..14) {
  $text = $kpiByRow[$r]
  $ws.Range("B$r").Value = $text
  $ws.Range("D$r").Value = $text
  $ws.Range("E$r").Value = $text
  $ws.Range("F$r").Value = $text
}

# Row 14 goes back to the sheet's standard row height (15) instead of the
# previous custom 13.8.
$ws.Rows.Item(14).RowHeight = 15

# Column widths were tightened (narrower) after the text changed.
$ws.Columns.Item(1).ColumnWidth = 15.8333333333333
$ws.Range("B1").EntireColumn.ColumnWidth = 60.1666666666667
$ws.Columns.Item(3).ColumnWidth = 17
$ws.Range("D1:F1").EntireColumn.ColumnWidth = 60.1666666666667
$ws.Range($ws.Cells.Item(1,7), $ws.Cells.Item(1,1025)).EntireColumn.ColumnWidth = 8.33333333333333

# Selection moves from a single cell (F14) to the whole data range (A1:F14).
$ws.Range("A1:F14").Select()
